$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 27.63641
$ws.Cells.Item(2, 8).Value = 82.90923000000001
$ws.Cells.Item(2, 9).Value = 0.1577716087636821
$ws.Cells.Item(2, 10).Value = 0.1653201523274948
$ws.Cells.Item(2, 13).Value = 2.906846333333333
$ws.Cells.Item(2, 14).Value = 8.720538999999999
$ws.Cells.Item(2, 15).Value = 0.005520525738044089
$ws.Cells.Item(2, 16).Value = 0.005624540846623205
$ws.Cells.Item(2, 17).Value = 80.33479707499666
$ws.Cells.Item(2, 18).Value = 723.0131736749699
$ws.Cells.Item(2, 19).Value = 0.0008709822269125293
$ws.Cells.Item(2, 20).Value = 0.0009298499495359648

$ws.Cells.Item(3, 7).Value = 27.63641
$ws.Cells.Item(3, 8).Value = 82.90923000000001
$ws.Cells.Item(3, 9).Value = 0.1577716087636821
$ws.Cells.Item(3, 10).Value = 0.1653201523274948
$ws.Cells.Item(3, 15).Value = 0.3528665483720876
$ws.Cells.Item(3, 16).Value = 0.3595150912979765
$ws.Cells.Item(3, 17).Value = 5134.920821521177
$ws.Cells.Item(3, 18).Value = 46214.2873936906
$ws.Cells.Item(3, 19).Value = 0.05567232301555192
$ws.Cells.Item(3, 20).Value = 0.05943508965741467

$ws.Cells.Item(4, 7).Value = 27.63641
$ws.Cells.Item(4, 8).Value = 82.90923000000001
$ws.Cells.Item(4, 9).Value = 0.1577716087636821
$ws.Cells.Item(4, 10).Value = 0.1653201523274948
$ws.Cells.Item(4, 13).Value = 137.0717086666666
$ws.Cells.Item(4, 14).Value = 411.2151259999999
$ws.Cells.Item(4, 15).Value = 0.2603191943704447
$ws.Cells.Item(4, 16).Value = 0.2652240042658267
$ws.Cells.Item(4, 17).Value = 3788.169940112553
$ws.Cells.Item(4, 18).Value = 34093.52946101298
$ws.Cells.Item(4, 19).Value = 0.04107097808789071
$ws.Cells.Item(4, 20).Value = 0.0438468727861346

$ws.Cells.Item(5, 7).Value = 27.63641
$ws.Cells.Item(5, 8).Value = 82.90923000000001
$ws.Cells.Item(5, 9).Value = 0.1577716087636821
$ws.Cells.Item(5, 10).Value = 0.1653201523274948
$ws.Cells.Item(5, 13).Value = 29.2127365
$ws.Cells.Item(5, 14).Value = 58.425473
$ws.Cells.Item(5, 15).Value = 0.05547925319534149
$ws.Cells.Item(5, 16).Value = 0.03768304451958546
$ws.Cells.Item(5, 17).Value = 807.335163135965
$ws.Cells.Item(5, 18).Value = 4844.010978815791
$ws.Cells.Item(5, 19).Value = 0.008753051029636676
$ws.Cells.Item(5, 20).Value = 0.006229766660141635

$ws.Cells.Item(6, 7).Value = 27.63641
$ws.Cells.Item(6, 8).Value = 82.90923000000001
$ws.Cells.Item(6, 9).Value = 0.1577716087636821
$ws.Cells.Item(6, 10).Value = 0.1653201523274948
$ws.Cells.Item(6, 13).Value = 171.5584106666666
$ws.Cells.Item(6, 14).Value = 514.6752319999999
$ws.Cells.Item(6, 15).Value = 0.3258144783240821
$ws.Cells.Item(6, 16).Value = 0.331953319069988
$ws.Cells.Item(6, 17).Value = 4741.258576132373
$ws.Cells.Item(6, 18).Value = 42671.32718519136
$ws.Cells.Item(6, 19).Value = 0.05140427440369026
$ws.Cells.Item(6, 20).Value = 0.0548785732742679

$ws.Cells.Item(7, 9).Value = 0.6207549685359464
$ws.Cells.Item(7, 10).Value = 0.650454836333234
$ws.Cells.Item(7, 13).Value = 2.906846333333333
$ws.Cells.Item(7, 14).Value = 8.720538999999999
$ws.Cells.Item(7, 15).Value = 0.005520525738044089
$ws.Cells.Item(7, 16).Value = 0.005624540846623205
$ws.Cells.Item(7, 17).Value = 316.078569657778
$ws.Cells.Item(7, 18).Value = 2844.707126920001
$ws.Cells.Item(7, 19).Value = 0.00342689378082144
$ws.Cells.Item(7, 20).Value = 0.003658509795839887

$ws.Cells.Item(8, 9).Value = 0.6207549685359464
$ws.Cells.Item(8, 10).Value = 0.650454836333234
$ws.Cells.Item(8, 15).Value = 0.3528665483720876
$ws.Cells.Item(8, 16).Value = 0.3595150912979765
$ws.Cells.Item(8, 19).Value = 0.2190436631321032
$ws.Cells.Item(8, 20).Value = 0.233848329869553

$ws.Cells.Item(9, 9).Value = 0.6207549685359464
$ws.Cells.Item(9, 10).Value = 0.650454836333234
$ws.Cells.Item(9, 13).Value = 137.0717086666666
$ws.Cells.Item(9, 14).Value = 411.2151259999999
$ws.Cells.Item(9, 15).Value = 0.2603191943704447
$ws.Cells.Item(9, 16).Value = 0.2652240042658267
$ws.Cells.Item(9, 17).Value = 14904.61642883805
$ws.Cells.Item(9, 18).Value = 134141.5478595425
$ws.Cells.Item(9, 19).Value = 0.1615944333107283
$ws.Cells.Item(9, 20).Value = 0.1725162362863732

$ws.Cells.Item(10, 9).Value = 0.6207549685359464
$ws.Cells.Item(10, 10).Value = 0.650454836333234
$ws.Cells.Item(10, 13).Value = 29.2127365
$ws.Cells.Item(10, 14).Value = 58.425473
$ws.Cells.Item(10, 15).Value = 0.05547925319534149
$ws.Cells.Item(10, 16).Value = 0.03768304451958546
$ws.Cells.Item(10, 17).Value = 3176.473370066769
$ws.Cells.Item(10, 18).Value = 19058.84022040061
$ws.Cells.Item(10, 19).Value = 0.03443902207167201
$ws.Cells.Item(10, 20).Value = 0.02451111855552493

$ws.Cells.Item(11, 9).Value = 0.6207549685359464
$ws.Cells.Item(11, 10).Value = 0.650454836333234
$ws.Cells.Item(11, 13).Value = 171.5584106666666
$ws.Cells.Item(11, 14).Value = 514.6752319999999
$ws.Cells.Item(11, 15).Value = 0.3258144783240821
$ws.Cells.Item(11, 16).Value = 0.331953319069988
$ws.Cells.Item(11, 17).Value = 18654.55921576006
$ws.Cells.Item(11, 18).Value = 167891.0329418406
$ws.Cells.Item(11, 19).Value = 0.2022509562406214
$ws.Cells.Item(11, 20).Value = 0.2159206418259429

$ws.Cells.Item(12, 7).Value = 6.867169333333333
$ws.Cells.Item(12, 8).Value = 20.601508
$ws.Cells.Item(12, 9).Value = 0.03920351280693195
$ws.Cells.Item(12, 10).Value = 0.04107919517207073
$ws.Cells.Item(12, 13).Value = 2.906846333333333
$ws.Cells.Item(12, 14).Value = 8.720538999999999
$ws.Cells.Item(12, 15).Value = 0.005520525738044089
$ws.Cells.Item(12, 16).Value = 0.005624540846623205
$ws.Cells.Item(12, 17).Value = 19.96180599697911
$ws.Cells.Item(12, 18).Value = 179.656253972812
$ws.Cells.Item(12, 19).Value = 0.0002164240014724089
$ws.Cells.Item(12, 20).Value = 0.0002310516111917186

$ws.Cells.Item(13, 7).Value = 6.867169333333333
$ws.Cells.Item(13, 8).Value = 20.601508
$ws.Cells.Item(13, 9).Value = 0.03920351280693195
$ws.Cells.Item(13, 10).Value = 0.04107919517207073
$ws.Cells.Item(13, 15).Value = 0.3528665483720876
$ws.Cells.Item(13, 16).Value = 0.3595150912979765
$ws.Cells.Item(13, 17).Value = 1275.93890793504
$ws.Cells.Item(13, 18).Value = 11483.45017141536
$ws.Cells.Item(13, 19).Value = 0.01383360824824301
$ws.Cells.Item(13, 20).Value = 0.01476859060273441

$ws.Cells.Item(14, 7).Value = 6.867169333333333
$ws.Cells.Item(14, 8).Value = 20.601508
$ws.Cells.Item(14, 9).Value = 0.03920351280693195
$ws.Cells.Item(14, 10).Value = 0.04107919517207073
$ws.Cells.Item(14, 13).Value = 137.0717086666666
$ws.Cells.Item(14, 14).Value = 411.2151259999999
$ws.Cells.Item(14, 15).Value = 0.2603191943704447
$ws.Cells.Item(14, 16).Value = 0.2652240042658267
$ws.Cells.Item(14, 17).Value = 941.294634223334
$ws.Cells.Item(14, 18).Value = 8471.651708010006
$ws.Cells.Item(14, 19).Value = 0.01020542687039193
$ws.Cells.Item(14, 20).Value = 0.01089518863555402

$ws.Cells.Item(15, 7).Value = 6.867169333333333
$ws.Cells.Item(15, 8).Value = 20.601508
$ws.Cells.Item(15, 9).Value = 0.03920351280693195
$ws.Cells.Item(15, 10).Value = 0.04107919517207073
$ws.Cells.Item(15, 13).Value = 29.2127365
$ws.Cells.Item(15, 14).Value = 58.425473
$ws.Cells.Item(15, 15).Value = 0.05547925319534149
$ws.Cells.Item(15, 16).Value = 0.03768304451958546
$ws.Cells.Item(15, 17).Value = 200.6088082355473
$ws.Cells.Item(15, 18).Value = 1203.652849413284
$ws.Cells.Item(15, 19).Value = 0.00217498161316259
$ws.Cells.Item(15, 20).Value = 0.001547989140497881

$ws.Cells.Item(16, 7).Value = 6.867169333333333
$ws.Cells.Item(16, 8).Value = 20.601508
$ws.Cells.Item(16, 9).Value = 0.03920351280693195
$ws.Cells.Item(16, 10).Value = 0.04107919517207073
$ws.Cells.Item(16, 13).Value = 171.5584106666666
$ws.Cells.Item(16, 14).Value = 514.6752319999999
$ws.Cells.Item(16, 15).Value = 0.3258144783240821
$ws.Cells.Item(16, 16).Value = 0.331953319069988
$ws.Cells.Item(16, 17).Value = 1178.120656605539
$ws.Cells.Item(16, 18).Value = 10603.08590944985
$ws.Cells.Item(16, 19).Value = 0.01277307207366201
$ws.Cells.Item(16, 20).Value = 0.01363637518209271

$ws.Cells.Item(17, 7).Value = 23.994483
$ws.Cells.Item(17, 8).Value = 47.988966
$ws.Cells.Item(17, 9).Value = 0.1369804610788022
$ws.Cells.Item(17, 10).Value = 0.09568950488575238
$ws.Cells.Item(17, 13).Value = 2.906846333333333
$ws.Cells.Item(17, 14).Value = 8.720538999999999
$ws.Cells.Item(17, 15).Value = 0.005520525738044089
$ws.Cells.Item(17, 16).Value = 0.005624540846623205
$ws.Cells.Item(17, 17).Value = 69.748274928779
$ws.Cells.Item(17, 18).Value = 418.489649572674
$ws.Cells.Item(17, 19).Value = 0.0007562041609946743
$ws.Cells.Item(17, 20).Value = 0.0005382095288230651

$ws.Cells.Item(18, 7).Value = 23.994483
$ws.Cells.Item(18, 8).Value = 47.988966
$ws.Cells.Item(18, 9).Value = 0.1369804610788022
$ws.Cells.Item(18, 10).Value = 0.09568950488575238
$ws.Cells.Item(18, 15).Value = 0.3528665483720876
$ws.Cells.Item(18, 16).Value = 0.3595150912979765
$ws.Cells.Item(18, 17).Value = 4458.240790259513
$ws.Cells.Item(18, 18).Value = 26749.44474155708
$ws.Cells.Item(18, 19).Value = 0.04833582249529404
$ws.Cells.Item(18, 20).Value = 0.03440182108525944

$ws.Cells.Item(19, 7).Value = 23.994483
$ws.Cells.Item(19, 8).Value = 47.988966
$ws.Cells.Item(19, 9).Value = 0.1369804610788022
$ws.Cells.Item(19, 10).Value = 0.09568950488575238
$ws.Cells.Item(19, 13).Value = 137.0717086666666
$ws.Cells.Item(19, 14).Value = 411.2151259999999
$ws.Cells.Item(19, 15).Value = 0.2603191943704447
$ws.Cells.Item(19, 16).Value = 0.2652240042658267
$ws.Cells.Item(19, 17).Value = 3288.964783383286
$ws.Cells.Item(19, 18).Value = 19733.78870029972
$ws.Cells.Item(19, 19).Value = 0.03565864327252585
$ws.Cells.Item(19, 20).Value = 0.02537915365201364

$ws.Cells.Item(20, 7).Value = 23.994483
$ws.Cells.Item(20, 8).Value = 47.988966
$ws.Cells.Item(20, 9).Value = 0.1369804610788022
$ws.Cells.Item(20, 10).Value = 0.09568950488575238
$ws.Cells.Item(20, 13).Value = 29.2127365
$ws.Cells.Item(20, 14).Value = 58.425473
$ws.Cells.Item(20, 15).Value = 0.05547925319534149
$ws.Cells.Item(20, 16).Value = 0.03768304451958546
$ws.Cells.Item(20, 17).Value = 700.9445093327296
$ws.Cells.Item(20, 18).Value = 2803.778037330918
$ws.Cells.Item(20, 19).Value = 0.007599573683005489
$ws.Cells.Item(20, 20).Value = 0.003605871872666897

$ws.Cells.Item(21, 7).Value = 23.994483
$ws.Cells.Item(21, 8).Value = 47.988966
$ws.Cells.Item(21, 9).Value = 0.1369804610788022
$ws.Cells.Item(21, 10).Value = 0.09568950488575238
$ws.Cells.Item(21, 13).Value = 171.5584106666666
$ws.Cells.Item(21, 14).Value = 514.6752319999999
$ws.Cells.Item(21, 15).Value = 0.3258144783240821
$ws.Cells.Item(21, 16).Value = 0.331953319069988
$ws.Cells.Item(21, 17).Value = 4116.455368248352
$ws.Cells.Item(21, 18).Value = 24698.73220949011
$ws.Cells.Item(21, 19).Value = 0.04463021746698219
$ws.Cells.Item(21, 20).Value = 0.03176444874698934

$ws.Cells.Item(22, 7).Value = 7.933225666666668
$ws.Cells.Item(22, 8).Value = 23.799677
$ws.Cells.Item(22, 9).Value = 0.04528944881463745
$ws.Cells.Item(22, 10).Value = 0.04745631128144808
$ws.Cells.Item(22, 13).Value = 2.906846333333333
$ws.Cells.Item(22, 14).Value = 8.720538999999999
$ws.Cells.Item(22, 15).Value = 0.005520525738044089
$ws.Cells.Item(22, 16).Value = 0.005624540846623205
$ws.Cells.Item(22, 17).Value = 23.06066794065589
$ws.Cells.Item(22, 18).Value = 207.546011465903
$ws.Cells.Item(22, 19).Value = 0.0002500215678430364
$ws.Cells.Item(22, 20).Value = 0.0002669199612325704

$ws.Cells.Item(23, 7).Value = 7.933225666666668
$ws.Cells.Item(23, 8).Value = 23.799677
$ws.Cells.Item(23, 9).Value = 0.04528944881463745
$ws.Cells.Item(23, 10).Value = 0.04745631128144808
$ws.Cells.Item(23, 15).Value = 0.3528665483720876
$ws.Cells.Item(23, 16).Value = 0.3595150912979765
$ws.Cells.Item(23, 17).Value = 1474.015100282305
$ws.Cells.Item(23, 18).Value = 13266.13590254074
$ws.Cells.Item(23, 19).Value = 0.01598113148089545
$ws.Cells.Item(23, 20).Value = 0.017061260083015

$ws.Cells.Item(24, 7).Value = 7.933225666666668
$ws.Cells.Item(24, 8).Value = 23.799677
$ws.Cells.Item(24, 9).Value = 0.04528944881463745
$ws.Cells.Item(24, 10).Value = 0.04745631128144808
$ws.Cells.Item(24, 13).Value = 137.0717086666666
$ws.Cells.Item(24, 14).Value = 411.2151259999999
$ws.Cells.Item(24, 15).Value = 0.2603191943704447
$ws.Cells.Item(24, 16).Value = 0.2652240042658267
$ws.Cells.Item(24, 17).Value = 1087.420797368256
$ws.Cells.Item(24, 18).Value = 9786.787176314301
$ws.Cells.Item(24, 19).Value = 0.01178971282890791
$ws.Cells.Item(24, 20).Value = 0.01258655290575119

$ws.Cells.Item(25, 7).Value = 7.933225666666668
$ws.Cells.Item(25, 8).Value = 23.799677
$ws.Cells.Item(25, 9).Value = 0.04528944881463745
$ws.Cells.Item(25, 10).Value = 0.04745631128144808
$ws.Cells.Item(25, 13).Value = 29.2127365
$ws.Cells.Item(25, 14).Value = 58.425473
$ws.Cells.Item(25, 15).Value = 0.05547925319534149
$ws.Cells.Item(25, 16).Value = 0.03768304451958546
$ws.Cells.Item(25, 17).Value = 231.7512309953702
$ws.Cells.Item(25, 18).Value = 1390.507385972221
$ws.Cells.Item(25, 19).Value = 0.002512624797864729
$ws.Cells.Item(25, 20).Value = 0.001788298290754114

$ws.Cells.Item(26, 7).Value = 7.933225666666668
$ws.Cells.Item(26, 8).Value = 23.799677
$ws.Cells.Item(26, 9).Value = 0.04528944881463745
$ws.Cells.Item(26, 10).Value = 0.04745631128144808
$ws.Cells.Item(26, 13).Value = 171.5584106666666
$ws.Cells.Item(26, 14).Value = 514.6752319999999
$ws.Cells.Item(26, 15).Value = 0.3258144783240821
$ws.Cells.Item(26, 16).Value = 0.331953319069988
$ws.Cells.Item(26, 17).Value = 1361.011586833341
$ws.Cells.Item(26, 18).Value = 12249.10428150006
$ws.Cells.Item(26, 19).Value = 0.01475595813912632
$ws.Cells.Item(26, 20).Value = 0.01575328004069521
